$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The report no longer needs first/last name or the created/updated
# timestamps - drop column A (Nombres) first so the surviving columns
# (Apellidos, Email, Creado el, Actualizado el) shift left, then drop
# the (now) trailing C:D date columns. Two sequential single-area
# deletes (rather than one multi-area delete) so both actually apply.
# This shifts the old Email column into B, preserving the existing
# (blank) email shared-string cell for row 2 instead of recreating it,
# since writing an empty string via .Value would simply clear the cell.
$ws.Range("A1:A1").EntireColumn.Delete()
$ws.Range("C1:D1").EntireColumn.Delete()

# New header for column A: the report now surfaces the record ID instead
# of the surname. Column B ("Email") already has the right header/values.
$ws.Range("A1").Value = "ID"

# Column A data becomes the record UUID for each row.
$ws.Range("A2").Value = "23c0d5f8-17c9-4c92-b743-f662561e554f"
$ws.Range("A3").Value = "19ca14e7-ace5-4d3b-8a6f-36ec768e9c5b"
$ws.Range("A4").Value = "20be17b5-065c-4d88-a3a9-6167ead560b3"
$ws.Range("A5").Value = "2d7f68de-5c96-42ab-86c2-5b6b5c2f2c0d"
$ws.Range("A6").Value = "34ca6cf3-1748-4c71-b6a2-3b4f2c48d3a9"
